$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph, before the "Get Ready to Embrace..." Heading2
#    paragraph.
# ---------------------------------------------------------------------
$titleHeading = $d.Paragraphs.Item(2)
$titleHeading.Range.InsertParagraphBefore()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.ParagraphFormat.Style = $d.Styles.Item("Normal")

# Insert the regular (non-bold) continuation text first...
$metaPara.Range.Text = ": Read our review of Celestial Beauty slot game from Skywind Group. Play for free with cascading mechanics and multiplier Wilds."

# ...then prepend "Meta description" at the very start of the paragraph...
$metaPara2 = $d.Paragraphs.Item(2)
$startPos = $metaPara2.Range.Start
$prependPoint = $d.Range($startPos, $startPos)
$prependPoint.InsertBefore("Meta description")

# ...and finally make only that prefix bold so the two parts remain separate runs.
$metaPara3 = $d.Paragraphs.Item(2)
$boldRange = $d.Range($metaPara3.Range.Start, $metaPara3.Range.Start + 16)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------
# 2. Remove the duplicate bold "Play Celestial Beauty Free Slot Game |
#    Review" paragraph near the end of the document.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping the italic run formatting intact.
# ---------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$descPara = $d.Paragraphs.Item($count2)
$newPromptText = "Please create a feature image fitting the game ""Celestial Beauty"" with the following specifications: - The image should be in a cartoon style - The image should feature a happy Maya warrior with glasses The image could feature the Maya warrior standing in a starry background, surrounded by celestial symbols, such as moons, suns, and stars. The warrior could be holding a staff or a sword, adorned with jewels and other precious stones. The image should be colorful and vibrant, capturing the excitement and energy of the game. The warrior should be drawn with a happy expression, conveying the enjoyment of playing the game. The glasses could be an important feature, adding a quirky and unique touch to the image. The Maya warrior could be portrayed as a powerful and brave figure, capturing the essence of the game's theme. Overall, the image should be visually appealing and engaging, appealing to players who enjoy the game's fun and light-hearted vibe."

$descTextRange = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descTextRange.Text = $newPromptText
